$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the short "urban"/"rural" labels (row 6 / row 7) with the fuller
# phrasing in all three languages (Kyrgyz / Russian / English columns).
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Update the view: drop the frozen/scrolled topLeftCell and move the
# active selection to B29.
$ws.Range("B29").Select()
